# Arbeitsjournal.xlsx update
# - "Andreas Lüscher" sheet: journal entries for the weeks since the last
#   recorded entry, plus the now-longer running total.
# - "Zusatz" sheet: SUMIF ranges widened to track the new data range.
# - "Arbeitsjournal" overview: picks up the new total via its formula.
# - Active sheet / selection bookkeeping to match where the author ended up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Andreas Lüscher" worksheet - append the new journal rows
# ---------------------------------------------------------------------------
$lu = $wb.Worksheets.Item("Andreas Lüscher")

# Existing row 49's topic was re-labelled.
$lu.Range("D49").Value = "Refactoring"

# New data rows 53-69 (previously blank template rows 53-57 plus brand new
# rows 58-69).
$rows = @(
    @{R=53; Date="02/07/2017"; C=6;  D="Writer";               E="Writer erweitert und neu strukturiert"},
    @{R=54; Date="02/09/2017"; C=8;  D="Unit Tests";           E="Unit Test mit Daten des Spektrometers erstellt und Writer verbessert"},
    @{R=55; Date="02/14/2017"; C=7;  D="UI";                   E="Settings Modals verbessert und Radiance View Controller geändert"},
    @{R=56; Date="02/15/2017"; C=7;  D="Prototyp 3";           E="Page Classes hinzugefügt für die einfachere verwaltung der Parameter"},
    @{R=57; Date="02/16/2017"; C=5;  D="Prototyp 3";           E="Spektras werden im Parent VC zwischengespeichert. Erst im Finish VC wird geschrieben."},
    @{R=58; Date="02/17/2017"; C=5;  D="UI & Controller";      E="RawSettings VC hinzugegfügt"},
    @{R=59; Date="02/20/2017"; C=4;  D="Writer";               E="Alle 3 Modi werden korrekt geschrieben"},
    @{R=60; Date="02/24/2017"; C=6;  D="Refactoring";          E="Code refactoring"},
    @{R=61; Date="02/25/2017"; C=4;  D="Calculations";         E="Neue Calculations Klasse hinzugefügt und Command Manager verbessert"},
    @{R=62; Date="02/26/2017"; C=4;  D="Writer";               E="Radiance wird korrekt geschrieben inklusive base lamp und fo files"},
    @{R=63; Date="02/28/2017"; C=9;  D="BackgroundFileManger"; E="Neue BackgroundFileMangaer Klasse hinzugefügt, um Daten bereits während dem Messen zu schreiben."},
    @{R=64; Date="03/01/2017"; C=12; D="Prototyp 3";           E="Verbesserte Calculations, Writes und Anzeige der Messresultate. Grosses refactoring der Messmethoden"},
    @{R=65; Date="03/02/2017"; C=7;  D="FileSelection";        E="Besuchte Pfade speichern und direkt dahinspringen, wenn FIleBrowser geöfnet wird."},
    @{R=66; Date="03/02/2017"; C=4;  D="File Managment";       E="Verbessertes FileManagement"},
    @{R=67; Date="03/03/2017"; C=5;  D="Instrument Control";   E="Instument Control Aktionen hinzugefügt"},
    @{R=68; Date="03/06/2017"; C=8;  D="Documentation";        E="Dokumentation erweitert"},
    @{R=69; Date="03/07/2017"; C=6;  D="Documentation";        E="Dokumentation erweitert"}
)

foreach ($row in $rows) {
    $r = $row.R
    $lu.Range("A$r").Formula = "=WEEKNUM(B$r)"
    $lu.Range("B$r").Value = $row.Date
    $lu.Range("C$r").Value = $row.C
    $lu.Range("D$r").Value = $row.D
    $lu.Range("E$r").Value = $row.E
}

# New trailing blank rows 70-73 (same pattern as the old blank rows used to
# have before the new data rows pushed them down).
for ($r = 70; $r -le 73; $r++) {
    $lu.Range("A$r").Formula = "=WEEKNUM(B$r)"
    $lu.Range("B$r").Value = ""
    $lu.Range("C$r").Value = ""
}
$lu.Range("C73").ClearContents()

# Move the summary / spacer block down by 16 rows (58->74 ... 62->78).
$lu.Range("B74").Value = $lu.Range("B58").Value
$lu.Range("C74").Formula = "=SUM(C4:C72)"
$lu.Range("B75").Value = $lu.Range("B59").Value
$lu.Range("B76").Value = $lu.Range("B60").Value
$lu.Range("B77").Value = $lu.Range("B61").Value
$lu.Range("B78").Value = $lu.Range("B62").Value

# Clear out the old summary row content that is now stale (row 58 already
# got overwritten with new journal content above, nothing further to clear).

$lu.Range("D73").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) "Zusatz" worksheet - widen the SUMIF ranges to the new data extent
# ---------------------------------------------------------------------------
$zs = $wb.Worksheets.Item("Zusatz")
for ($r = 2; $r -le 26; $r++) {
    $zs.Range("F$r").Formula = "=SUMIF('Andreas Lüscher'!`$A`$4:`$C`$72,A$r,'Andreas Lüscher'!`$C`$4:`$C`$72)"
}

# ---------------------------------------------------------------------------
# 3) "Arbeitsjournal" overview - formula now points at the new total cell
# ---------------------------------------------------------------------------
$aj = $wb.Worksheets.Item("Arbeitsjournal")
$aj.Range("J5").Formula = "='Andreas Lüscher'!C74"
$aj.Range("M26").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Active sheet bookkeeping - author ended up on "Andreas Lüscher"
# ---------------------------------------------------------------------------
$lu.Activate()
$lu.Range("D73").Select() | Out-Null
